$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 351; existing rows 351-415 shift down to 352-416
$ws.Rows("351:351").Insert()

# Populate new row 351 with the new record's data
$ws.Range("A351").Value = 4
$ws.Range("B351").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C351").Value = "Los Lagos"
$ws.Range("D351").Value = 44995
$ws.Range("D351").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E351").Value = 10
$ws.Range("F351").Value = 100112040
$ws.Range("G351").Value = "Cilantro"
$ws.Range("H351").Value = "Sin especificar"
$ws.Range("I351").Value = "Primera"
$ws.Range("J351").Value = 140
$ws.Range("K351").Value = 12000
$ws.Range("L351").Value = 14000
$ws.Range("M351").Value = 13000
$ws.Range("N351").Value = "$/caja 36 atados"
$ws.Range("O351").Value = "Región Metropolitana"
$ws.Range("P351").Value = 361
$ws.Range("Q351").Value = 36
$ws.Range("R351").Value = "Hortaliza"
